$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Cognitive Failures
$ws.Range("B2").Value = 89
$ws.Range("C2").Value = 2.6
$ws.Range("D2").Value = 0.9

# Row 3: Cognitive Load
$ws.Range("B3").Value = 89
$ws.Range("C3").Value = 2.04
$ws.Range("D3").Value = 0.8100000000000001

# Row 4: Performance Expectancy
$ws.Range("B4").Value = 89
$ws.Range("C4").Value = 2.6

# Row 5: Effort Expectancy
$ws.Range("B5").Value = 89
$ws.Range("C5").Value = 3.73
$ws.Range("D5").Value = 1.16

# Row 6: Facilitating Conditions
$ws.Range("B6").Value = 89
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 1.02

# Row 7: Trust in AI (core)
$ws.Range("B7").Value = 89
$ws.Range("C7").Value = 2.79
$ws.Range("D7").Value = 0.9399999999999999

# Row 8: Checking Behaviour
$ws.Range("B8").Value = 89
$ws.Range("C8").Value = 3.62
$ws.Range("D8").Value = 1.48
$ws.Range("F8").Value = 3
